$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final layout (rows 1-27), columns A (label), B (06-01-2023), C (13-01-2023)
# Row 1 is the header row with the two report dates in B1 / C1 (A1 stays empty).
# Rows 2-25 are the individual funds (alphabetical, unchanged from before except for
# the two summary rows 'total' and 'avg' which move from rows 2-3 down to rows 26-27).
# A new column C is added with the values as of 13-01-2023.

$labels = @(
  "1810 Renta variable",
  "1822 Raices Valores Negociables",
  "Alpha Acciones",
  "Alpha Mega",
  "Alpha planeam equil",
  "Argenfunds",
  "Balanz",
  "Consultatio Renta Variable",
  "FBA Acciones Argentinas",
  "FBA Calificado",
  "Fima Acciones",
  "Fima PB Acciones",
  "Gainvest Renta Variable",
  "Goal Acciones Argentinas",
  "Goal acciones plus",
  "HF Acciones Argentinas",
  "IAM Renta Variable",
  "IEB Value",
  "Lombardi",
  "Pellegrini Acciones",
  "Pionero Acciones",
  "Premier Renta Variable",
  "Rofex 20 Renta Variable",
  "Supefondo RV",
  "avg",
  "total"
)

$colB = @(
  197337.63,
  221699.36,
  127993.49,
  217335.35,
  8468.639999999999,
  19056.78,
  156144.15,
  4387.69,
  142260.66,
  138362.63,
  555525.9,
  155940.79,
  150378.55,
  5156.3,
  9287.01,
  1034.03,
  61565.16,
  6575.87,
  29486.84,
  78133.67,
  34736.42,
  23265.9,
  136656.94,
  691621.47,
  132183.8,
  3172411.23
)

$colC = @(
  196831.15,
  196327.84,
  128192.98,
  218178.8,
  8405.610000000001,
  19057.86,
  155761.69,
  4514.12,
  141683.49,
  138837.38,
  555676.76,
  154529.9,
  150329.28,
  5129.72,
  9324.84,
  1067.42,
  64384.36,
  6617.12,
  29464.8,
  78148.84,
  35021.11,
  23466.91,
  136740.05,
  687495.29,
  131049.47,
  3145187.32
)

# Header row: B1 already holds "06-01-2023" and is left untouched so Excel
# does not reinterpret/reformat it. Only the new C1 date is added, copying
# B1's formatting (bold, centered, bordered) without touching B1 itself.
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $labels.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $labels[$i]
  $ws.Cells.Item($row, 2).Value = $colB[$i]
  $ws.Cells.Item($row, 3).Value = $colC[$i]
}
